$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text tweaks (rich-text runs: replace only the trailing digits/dates) ---
# A8: "Volume 30   Number  5" -> "...6"  (the final run, "5", becomes "6")
$ws.Range("A8").Characters(21, 1).Text = "6"

# C9: "Report Covering the Week  1/30/2023  Through  2/5/2023"
#     -> "...2/6/2023  Through  2/12/2023"
$ws.Range("C9").Characters(27, 9).Text = "2/6/2023"
$ws.Range("C9").Characters(46, 8).Text = "2/12/2023"

# --- Row 19 (Gr. Larceny) ---
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 3
$ws.Range("K19").Value = 0
$ws.Range("N19").Value = -62.5

# --- Row 20 (G.L.A.) : N20 goes from blank placeholder text to a real percentage ---
$ws.Range("N20").Value = -100
$ws.Range("N20").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Row 21 (TOTAL) ---
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = -33.333333333333
$ws.Range("J21").Value = 5
$ws.Range("K21").Value = -20
$ws.Range("N21").Value = -87.5

# --- Row 24 (Petit Larceny) : C24 goes from blank placeholder text to a real count ---
$ws.Range("C24").Value = 2
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("F24").Value = 5
$ws.Range("I24").Value = 5
$ws.Range("K24").Value = 400
$ws.Range("L24").Value = 150
$ws.Range("M24").Value = 400

# --- Row 25 (Misd. Assault) : C25 goes from blank placeholder text to a real count ---
$ws.Range("C25").Value = 2
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("F25").Value = 2
$ws.Range("I25").Value = 3
$ws.Range("K25").Value = 50
$ws.Range("L25").Value = 200
$ws.Range("M25").Value = 0
